# Update the "想去人数" (interested-count) figures in column F of the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets to the refreshed
# values from the latest data pull.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 2960
$ws1.Range("F5").Value  = 211
$ws1.Range("F6").Value  = 109
$ws1.Range("F8").Value  = 1632
$ws1.Range("F9").Value  = 1605
$ws1.Range("F10").Value = 52
$ws1.Range("F11").Value = 349
$ws1.Range("F18").Value = 226
$ws1.Range("F19").Value = 20
$ws1.Range("F20").Value = 11
$ws1.Range("F21").Value = 35
$ws1.Range("F22").Value = 4
$ws1.Range("F23").Value = 350
$ws1.Range("F24").Value = 136
$ws1.Range("F25").Value = 92
$ws1.Range("F27").Value = 1972
$ws1.Range("F29").Value = 451
$ws1.Range("F30").Value = 12
$ws1.Range("F31").Value = 162
$ws1.Range("F32").Value = 571
$ws1.Range("F33").Value = 220
$ws1.Range("F34").Value = 328
$ws1.Range("F36").Value = 484

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 2960
$ws4.Range("F5").Value  = 211
$ws4.Range("F6").Value  = 109
$ws4.Range("F8").Value  = 1632
$ws4.Range("F9").Value  = 1605
$ws4.Range("F10").Value = 52
$ws4.Range("F11").Value = 349
$ws4.Range("F18").Value = 226
$ws4.Range("F19").Value = 20
$ws4.Range("F20").Value = 11
$ws4.Range("F21").Value = 35
$ws4.Range("F22").Value = 4
$ws4.Range("F23").Value = 350
$ws4.Range("F24").Value = 136
$ws4.Range("F25").Value = 92
$ws4.Range("F27").Value = 1973
$ws4.Range("F29").Value = 451
$ws4.Range("F30").Value = 12
$ws4.Range("F31").Value = 162
$ws4.Range("F32").Value = 571
$ws4.Range("F33").Value = 220
$ws4.Range("F34").Value = 328
$ws4.Range("F36").Value = 484
